$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value (45192 = 2023-09-23)
# that needs to be bumped to 45202 (2023-10-03) for every data row (2..484).
$ws.Range("C2:C484").Value = 45202
